# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the three
# new header cells so they pick up the same bold/border/centered style used
# by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-49) gets the team's season record repeated across it.
$firstDataRow = 2
$lastDataRow = 49
$winsCol = 30   # AD
$lossesCol = 31 # AE
$tiesCol = 32   # AF

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, $winsCol).Value = 95
    $ws.Cells.Item($r, $lossesCol).Value = 65
    $ws.Cells.Item($r, $tiesCol).Value = 1
}
